# Zeiterfassung.xlsx - "Plugin Temperatur, Static, Navi - Ready" update
#
# Adds four new time-tracking entries (rows 15-18) to the "Zeitliste"
# worksheet for the work done on the Static/Navi plugins, carrying over
# the same date / number formatting used by the existing entries above
# them, and leaves the cursor on the next empty row (E19) as the author
# did after finishing their data entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Zeitliste")

# --- Row 15: 2013-11-22, Daniel Hörmann, Development, 5h, Plugin Navi ---
$ws.Range("A15").Value = 41600
$ws.Range("B15").Value = "Daniel Hörmann"
$ws.Range("C15").Value = "Development"
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = "Plugin Navi"

# --- Row 16: 2013-11-23, Daniel Hörmann, Development, 6h, Plugin Static flush fix ---
$ws.Range("A16").Value = 41601
$ws.Range("B16").Value = "Daniel Hörmann"
$ws.Range("C16").Value = "Development"
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = "Plugin Static - Probleme bei der Bildübertragung - Flush "

# --- Row 17: 2013-11-24, Daniel Hörmann, Development, 5h, Plugin Static + Navi ---
$ws.Range("A17").Value = 41602
$ws.Range("B17").Value = "Daniel Hörmann"
$ws.Range("C17").Value = "Development"
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = "Plugin Static, Plugin Navi"

# --- Row 18: 2013-11-25, Daniel Hörmann, Development, 3h, Plugin Navi multiuser ---
$ws.Range("A18").Value = 41603
$ws.Range("B18").Value = "Daniel Hörmann"
$ws.Range("C18").Value = "Development"
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = "Plugin Navi - Multiuserfähigkeit"

# Carry the date number format from the last filled-in entry (row 12) down
# onto the new dates so A15:A18 render the same way as A5:A12.
$ws.Range("A12").Copy()
$ws.Range("A15:A18").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# The "Multiuserfähigkeit" note ended up with an (accidental) time format
# applied to it, matching the source edit.
$ws.Range("E18").NumberFormat = "h:mm"

# Leave the selection on the next empty row, where the author's cursor was
# after typing the new entries.
$ws.Range("E19").Select() | Out-Null
